$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("JSS 3E")

# Fix existing row 7: "Admission No" (C7) should be stored as a real number (33)
$ws.Cells.Item(7, 3).Value = 33

# Append the new submission as row 8
$ws.Cells.Item(8, 1).Value = "2026-02-08 19:38:02"
$ws.Cells.Item(8, 2).Value = "RAHAMA MOHAMMED YUSUF"

# Admission No for this submission stays textual ("42"), matching the source export
$ws.Cells.Item(8, 3).NumberFormat = "@"
$ws.Cells.Item(8, 3).Value = "42"
$ws.Cells.Item(8, 3).Style = "Normal"

$ws.Cells.Item(8, 4).Value = 8
